$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "continent_id" header (column A) to "continent_code"
$ws.Range("A1").Value = "continent_code"

# The sheet used to carry a bunch of trailing, value-less columns (F:J)
# that only held stray formatting. Drop them completely.
$ws.Range("F1:J2").Clear()

# Column E ("sort_order") only has a header - row 2 never used it.
# Remove the stray formatted-but-empty cell at E2.
$ws.Range("E2").Clear()

# Re-apply E1's value so the cell regains the same formatting style as
# the rest of the header row (A1:D1) instead of the old, now-unused style.
$sortOrderHeader = $ws.Range("E1").Value2
$ws.Range("E1").Clear()
$ws.Range("E1").Value = $sortOrderHeader

# Reset the active selection back to A1 (top-left), matching the
# refreshed import sheet.
[void]$ws.Range("A1").Select()
